# Build site at 2022-09-26 16:07:08 UTC
# Rearranges the "Objetivos / Programa / Avaliacao" block (rows 10-22) and
# removes the old "Requisitos:" label row (old row 23), which shifts the
# whole requirement list up by one and drops the final (now-empty) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos: value replaced with the supervisor's name ---
$ws.Range("B10").Value = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Range("C10").Value = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Rows(10).RowHeight = 60

# --- Row 11: Objectives: (unchanged, still empty value) ---
$ws.Rows(11).RowHeight = 60

# --- Row 12: Docentes responsaveis: (unchanged, still empty value) ---

# --- Row 13: now "Programa resumido:" / "Semestral" ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows(13).RowHeight = 60

# --- Row 14: now "Short syllabus:" with no value ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Rows(14).RowHeight = 60

# --- Row 15: now "Programa:" / "01/01/2022" ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2022"
$ws.Range("C15").Value = "01/01/2022"
$ws.Rows(15).RowHeight = 120

# --- Row 16: now "Syllabus:" with no value ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Rows(16).RowHeight = 120

# --- Row 17: now "Avaliacao:" with no value/height ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows(17).RowHeight = 15

# --- Row 18: now "Metodo:" / supervisor's name ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Range("C18").Value = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Rows(18).RowHeight = 60

# --- Row 19: now "Criterio:" / supervision text ---
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
$ws.Range("C19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
$ws.Rows(19).RowHeight = 60

# --- Row 20: now "Norma de recuperacao:" / final grade text ---
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A nota final será baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio."
$ws.Range("C20").Value = "A nota final será baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio."
$ws.Rows(20).RowHeight = 60

# --- Row 21: now "Bibliografia:" / no-recovery text ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Devido às características da disciplina, não será oferecida recuperação."
$ws.Range("C21").Value = "Devido às características da disciplina, não será oferecida recuperação."
$ws.Rows(21).RowHeight = 120

# --- Row 22: now "Requisitos:" with no value/height (was Bibliografia/"Não há.") ---
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Rows(22).RowHeight = 15

# --- Old row 23 ("Requisitos:" label, now duplicated at row 22) is removed.
# This shifts the whole requirement list (old rows 24-49) up by one row and
# drops the final, now-blank row off the bottom (49 -> 48 total rows).
$ws.Rows(23).Delete()
